# Add pair tests (1.3.6 and 1.3.7) for claim functionality:
# mark rows 23 and 24 (column D, "Автоматизировано:") with "V"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D23").Value = "V"
$ws.Range("D24").Value = "V"

# Copy the style used by the other "Автоматизировано" cells (e.g. D22) so the
# new cells match formatting.
$ws.Range("D22").Copy() | Out-Null
$ws.Range("D23:D24").PasteSpecial(-4122) | Out-Null # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("D27").Select() | Out-Null
